$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Bugs and errors")
$ws2 = $wb.Worksheets.Item("Implemented Features")

# --- "Implemented Features" sheet content updates ---
# Row 5: "Added python tracking file" -> "Added python changelog excel file"
$ws2.Range("A5").Value = "Added python changelog excel file"
# Row 7: "Var / mean of inning" -> "f.calc_variability_seg_M_joint"
$ws2.Range("A7").Value = "f.calc_variability_seg_M_joint"
# Row 9: "Rotation of new optitrack dataset" -> "f.orient_markers", with a new comment
$ws2.Range("A9").Value = "f.orient_markers"
$ws2.Range("D9").Value = "Orients markers to match old data set"

# Widen column D on "Implemented Features" (target stored width 43.42578125;
# engine snaps column widths to coarse increments, closest achievable is ~43.5)
$ws2.Columns.Item(4).ColumnWidth = 42.666666667

# --- Selection / active sheet changes ---
# "Implemented Features" loses the tab selection / its own selection moves to B19
$ws2.Range("B19").Select() | Out-Null
# "Bugs and errors" becomes the active/selected tab with selection at C15
$ws1.Select() | Out-Null
$ws1.Range("C15").Select() | Out-Null
